# Auto-generated edit script: updates market-price derived cells
# (currentAveragePrice / LevePrice / LeveProfit columns H-N) across
# all 8 leve sheets, matching the scheduled market-data refresh diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 599.5
$ws.Range("I20").Value = 599.5
$ws.Range("K20").Value = 599.5
$ws.Range("M20").Value = -369.5
$ws.Range("H28").Value = 12825770
$ws.Range("I28").Value = 18525414
$ws.Range("J28").Value = 1571.125
$ws.Range("K28").Value = 18525414
$ws.Range("L28").Value = 1571.125
$ws.Range("M28").Value = -18524929
$ws.Range("N28").Value = -2541.125
$ws.Range("H35").Value = 599.5
$ws.Range("I35").Value = 599.5
$ws.Range("K35").Value = 599.5
$ws.Range("M35").Value = -220.5
$ws.Range("H41").Value = 1209.2667
$ws.Range("I41").Value = 1875.3334
$ws.Range("J41").Value = 765.2222
$ws.Range("K41").Value = 1875.3334
$ws.Range("L41").Value = 765.2222
$ws.Range("M41").Value = -1435.3334
$ws.Range("N41").Value = -1645.2222
$ws.Range("H70").Value = 1701
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 2051.5
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 6154.5
$ws.Range("M70").Value = -2730
$ws.Range("N70").Value = -6694.5
$ws.Range("H73").Value = 1701
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 2051.5
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 6154.5
$ws.Range("M73").Value = -2064
$ws.Range("N73").Value = -8026.5
$ws.Range("H80").Value = 2249.1177
$ws.Range("I80").Value = 572
$ws.Range("J80").Value = 4645
$ws.Range("K80").Value = 1716
$ws.Range("L80").Value = 13935
$ws.Range("M80").Value = -718
$ws.Range("N80").Value = -15931
$ws.Range("H83").Value = 2249.1177
$ws.Range("I83").Value = 572
$ws.Range("J83").Value = 4645
$ws.Range("K83").Value = 5148
$ws.Range("L83").Value = 41805
$ws.Range("M83").Value = -156
$ws.Range("N83").Value = -51789
$ws.Range("H86").Value = 2089.682
$ws.Range("I86").Value = 1807.6875
$ws.Range("J86").Value = 2841.6667
$ws.Range("K86").Value = 1807.6875
$ws.Range("L86").Value = 2841.6667
$ws.Range("M86").Value = -684.6875
$ws.Range("N86").Value = -5087.6667
$ws.Range("H88").Value = 4045578.5
$ws.Range("I88").Value = 1536.1428
$ws.Range("J88").Value = 5177910.5
$ws.Range("K88").Value = 1536.1428
$ws.Range("L88").Value = 5177910.5
$ws.Range("M88").Value = -1130.1428
$ws.Range("N88").Value = -5178722.5
$ws.Range("H89").Value = 2089.682
$ws.Range("I89").Value = 1807.6875
$ws.Range("J89").Value = 2841.6667
$ws.Range("K89").Value = 9038.4375
$ws.Range("L89").Value = 14208.3335
$ws.Range("M89").Value = -3422.4375
$ws.Range("N89").Value = -25440.3335
$ws.Range("H91").Value = 4045578.5
$ws.Range("I91").Value = 1536.1428
$ws.Range("J91").Value = 5177910.5
$ws.Range("K91").Value = 1536.1428
$ws.Range("L91").Value = 5177910.5
$ws.Range("M91").Value = -132.1428000000001
$ws.Range("N91").Value = -5180718.5
$ws.Range("H116").Value = 10683
$ws.Range("I116").Value = 26251.25
$ws.Range("J116").Value = 2898.875
$ws.Range("K116").Value = 26251.25
$ws.Range("L116").Value = 2898.875
$ws.Range("M116").Value = -22809.25
$ws.Range("N116").Value = -9782.875
$ws.Range("H132").Value = 29413910
$ws.Range("I132").Value = 41667600
$ws.Range("J132").Value = 5055.2
$ws.Range("K132").Value = 125002800
$ws.Range("L132").Value = 15165.6
$ws.Range("M132").Value = -125000270
$ws.Range("N132").Value = -20225.6
$ws.Range("H133").Value = 59625
$ws.Range("J133").Value = 59625
$ws.Range("L133").Value = 59625
$ws.Range("N133").Value = -69745

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 142863260
$ws.Range("I37").Value = 333336670
$ws.Range("J37").Value = 8200.75
$ws.Range("K37").Value = 333336670
$ws.Range("L37").Value = 8200.75
$ws.Range("M37").Value = -333336397
$ws.Range("N37").Value = -8746.75
$ws.Range("H61").Value = 997.2632
$ws.Range("I61").Value = 644.0476
$ws.Range("K61").Value = 644.0476
$ws.Range("M61").Value = -432.0476
$ws.Range("H119").Value = 30887
$ws.Range("J119").Value = 30887
$ws.Range("L119").Value = 30887
$ws.Range("N119").Value = -40563
$ws.Range("H136").Value = 997.2632
$ws.Range("I136").Value = 644.0476
$ws.Range("K136").Value = 1932.1428
$ws.Range("M136").Value = 617.8571999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 17713.834
$ws.Range("I82").Value = 5800
$ws.Range("J82").Value = 21685.111
$ws.Range("K82").Value = 5800
$ws.Range("L82").Value = 21685.111
$ws.Range("M82").Value = -5417
$ws.Range("N82").Value = -22451.111
$ws.Range("H85").Value = 17713.834
$ws.Range("I85").Value = 5800
$ws.Range("J85").Value = 21685.111
$ws.Range("K85").Value = 5800
$ws.Range("L85").Value = 21685.111
$ws.Range("M85").Value = -4474
$ws.Range("N85").Value = -24337.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1875.762
$ws.Range("I16").Value = 1870.9286
$ws.Range("J16").Value = 1885.4286
$ws.Range("K16").Value = 1870.9286
$ws.Range("L16").Value = 1885.4286
$ws.Range("M16").Value = -1583.9286
$ws.Range("N16").Value = -2459.4286
$ws.Range("H22").Value = 565
$ws.Range("I22").Value = 537.7778
$ws.Range("J22").Value = 600
$ws.Range("K22").Value = 537.7778
$ws.Range("L22").Value = 600
$ws.Range("M22").Value = -187.7778
$ws.Range("N22").Value = -1300
$ws.Range("H60").Value = 8450.75
$ws.Range("J60").Value = 8450.75
$ws.Range("L60").Value = 8450.75
$ws.Range("N60").Value = -9472.75
$ws.Range("H94").Value = 4673.25
$ws.Range("I94").Value = 458.66666
$ws.Range("J94").Value = 7202
$ws.Range("K94").Value = 458.66666
$ws.Range("L94").Value = 7202
$ws.Range("M94").Value = -7.666659999999979
$ws.Range("N94").Value = -8104
$ws.Range("H107").Value = 448.5
$ws.Range("I107").Value = 374.42105
$ws.Range("J107").Value = 730
$ws.Range("K107").Value = 374.42105
$ws.Range("L107").Value = 730
$ws.Range("M107").Value = 1545.57895
$ws.Range("N107").Value = -4570
$ws.Range("H113").Value = 1875.762
$ws.Range("I113").Value = 1870.9286
$ws.Range("J113").Value = 1885.4286
$ws.Range("K113").Value = 1870.9286
$ws.Range("L113").Value = 1885.4286
$ws.Range("M113").Value = 299.0714
$ws.Range("N113").Value = -6225.4286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 972.7059
$ws.Range("I113").Value = 772.7143
$ws.Range("J113").Value = 1004.5227
$ws.Range("K113").Value = 2318.1429
$ws.Range("L113").Value = 3013.5681
$ws.Range("M113").Value = -148.1428999999998
$ws.Range("N113").Value = -7353.5681
$ws.Range("H134").Value = 5381.6763
$ws.Range("J134").Value = 4792.857
$ws.Range("L134").Value = 14378.571
$ws.Range("N134").Value = -24518.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 29473.334
$ws.Range("J119").Value = 29473.334
$ws.Range("L119").Value = 29473.334
$ws.Range("N119").Value = -39149.334
$ws.Range("H132").Value = 3295.822
$ws.Range("I132").Value = 3694.5088
$ws.Range("J132").Value = 1875.5
$ws.Range("K132").Value = 11083.5264
$ws.Range("L132").Value = 5626.5
$ws.Range("M132").Value = -8553.526400000001
$ws.Range("N132").Value = -10686.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3497.889
$ws.Range("I81").Value = 2615.7778
$ws.Range("J81").Value = 4380
$ws.Range("K81").Value = 5231.5556
$ws.Range("L81").Value = 8760
$ws.Range("M81").Value = -4170.5556
$ws.Range("N81").Value = -10882
$ws.Range("H84").Value = 3497.889
$ws.Range("I84").Value = 2615.7778
$ws.Range("J84").Value = 4380
$ws.Range("K84").Value = 26157.778
$ws.Range("L84").Value = 43800
$ws.Range("M84").Value = -20853.778
$ws.Range("N84").Value = -54408
$ws.Range("H119").Value = 21810
$ws.Range("J119").Value = 21810
$ws.Range("L119").Value = 21810
$ws.Range("N119").Value = -31486
$ws.Range("H135").Value = 57016.668
$ws.Range("J135").Value = 57016.668
$ws.Range("L135").Value = 57016.668
$ws.Range("N135").Value = -67156.66800000001
